$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 41.43823833333334
$ws.Range("H2").Value = 124.314715
$ws.Range("I2").Value = 0.981992391336623
$ws.Range("J2").Value = 0.9819923913366232
$ws.Range("M2").Value = 31.618405
$ws.Range("N2").Value = 94.855215
$ws.Range("O2").Value = 0.8578613706944929
$ws.Range("P2").Value = 0.8578613706944929
$ws.Range("Q2").Value = 1310.211002109859
$ws.Range("R2").Value = 11791.89901898873
$ws.Range("S2").Value = 0.8424133388435983
$ws.Range("T2").Value = 0.8424133388435984

# Row 3
$ws.Range("G3").Value = 41.43823833333334
$ws.Range("H3").Value = 124.314715
$ws.Range("I3").Value = 0.981992391336623
$ws.Range("J3").Value = 0.9819923913366232
$ws.Range("O3").Value = 0.08747555172986397
$ws.Range("P3").Value = 0.08747555172986396
$ws.Range("Q3").Value = 133.6013419036606
$ws.Range("R3").Value = 1202.412077132945
$ws.Range("S3").Value = 0.08590032622669959
$ws.Range("T3").Value = 0.08590032622669959

# Row 4
$ws.Range("G4").Value = 41.43823833333334
$ws.Range("H4").Value = 124.314715
$ws.Range("I4").Value = 0.981992391336623
$ws.Range("J4").Value = 0.9819923913366232
$ws.Range("M4").Value = 2.014730333333334
$ws.Range("N4").Value = 6.044191000000001
$ws.Range("O4").Value = 0.05466307757564324
$ws.Range("P4").Value = 0.05466307757564324
$ws.Range("Q4").Value = 83.48687573006279
$ws.Range("R4").Value = 751.3818815705652
$ws.Range("S4").Value = 0.05367872626632524
$ws.Range("T4").Value = 0.05367872626632525

# Row 5
$ws.Range("I5").Value = 0.006845967574057415
$ws.Range("J5").Value = 0.006845967574057417
$ws.Range("M5").Value = 31.618405
$ws.Range("N5").Value = 94.855215
$ws.Range("O5").Value = 0.8578613706944929
$ws.Range("P5").Value = 0.8578613706944929
$ws.Range("Q5").Value = 9.134146165235
$ws.Range("R5").Value = 82.207315487115
$ws.Range("S5").Value = 0.005872891126810946
$ws.Range("T5").Value = 0.005872891126810948

# Row 6
$ws.Range("I6").Value = 0.006845967574057415
$ws.Range("J6").Value = 0.006845967574057417
$ws.Range("O6").Value = 0.08747555172986397
$ws.Range("P6").Value = 0.08747555172986396
$ws.Range("S6").Value = 0.0005988547906654308
$ws.Range("T6").Value = 0.0005988547906654308

# Row 7
$ws.Range("I7").Value = 0.006845967574057415
$ws.Range("J7").Value = 0.006845967574057417
$ws.Range("M7").Value = 2.014730333333334
$ws.Range("N7").Value = 6.044191000000001
$ws.Range("O7").Value = 0.05466307757564324
$ws.Range("P7").Value = 0.05466307757564324
$ws.Range("Q7").Value = 0.5820294018056668
$ws.Range("S7").Value = 0.0003742216565810386
$ws.Range("T7").Value = 0.0003742216565810387

# Row 8
$ws.Range("I8").Value = 0.01116164108931947
$ws.Range("J8").Value = 0.01116164108931947
$ws.Range("M8").Value = 31.618405
$ws.Range("N8").Value = 94.855215
$ws.Range("O8").Value = 0.8578613706944929
$ws.Range("P8").Value = 0.8578613706944929
$ws.Range("Q8").Value = 14.89227929446833
$ws.Range("R8").Value = 134.030513650215
$ws.Range("S8").Value = 0.009575140724083572
$ws.Range("T8").Value = 0.009575140724083575

# Row 9
$ws.Range("I9").Value = 0.01116164108931947
$ws.Range("J9").Value = 0.01116164108931947
$ws.Range("O9").Value = 0.08747555172986397
$ws.Range("P9").Value = 0.08747555172986396
$ws.Range("S9").Value = 0.0009763707124989404
$ws.Range("T9").Value = 0.0009763707124989408

# Row 10
$ws.Range("I10").Value = 0.01116164108931947
$ws.Range("J10").Value = 0.01116164108931947
$ws.Range("M10").Value = 2.014730333333334
$ws.Range("N10").Value = 6.044191000000001
$ws.Range("O10").Value = 0.05466307757564324
$ws.Range("P10").Value = 0.05466307757564324
$ws.Range("Q10").Value = 0.9489386585767778
$ws.Range("S10").Value = 0.0006101296527369572
$ws.Range("T10").Value = 0.0006101296527369575
